$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Hassan"
$ws.Range("C3").Value = "SP-001"
$ws.Range("D3").Value = "d"
$ws.Range("E3").Value = 200
$ws.Range("F3").Value = 30
$ws.Range("G3").Value = 12
$ws.Range("H3").Value = "piece"
$ws.Range("I3").Value = 300
$ws.Range("J3").Value = 18
$ws.Range("K3").Value = $true
